# Update report data through 2021-12-08 (rows 386-464)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @(44460, 1, 6, 91.37983551629607),
    @(44461, 0, 6, 91.37983551629607),
    @(44462, 1, 6, 91.37983551629607),
    @(44463, 1, 7, 106.6098081023454),
    @(44464, 0, 4, 60.91989034419738),
    @(44465, 1, 4, 60.91989034419738),
    @(44466, 0, 4, 60.91989034419738),
    @(44467, 0, 3, 45.68991775814803),
    @(44468, 0, 3, 45.68991775814803),
    @(44469, 0, 2, 30.45994517209869),
    @(44470, 0, 1, 15.22997258604935),
    @(44471, 0, 1, 15.22997258604935),
    @(44472, 0, 0, 0),
    @(44473, 2, 2, 30.45994517209869),
    @(44474, 0, 2, 30.45994517209869),
    @(44475, 0, 2, 30.45994517209869),
    @(44476, 0, 2, 30.45994517209869),
    @(44477, 0, 2, 30.45994517209869),
    @(44478, 0, 2, 30.45994517209869),
    @(44479, 0, 2, 30.45994517209869),
    @(44480, 0, 0, 0),
    @(44481, 0, 0, 0),
    @(44482, 0, 0, 0),
    @(44483, 0, 0, 0),
    @(44484, 0, 0, 0),
    @(44485, 0, 0, 0),
    @(44486, 0, 0, 0),
    @(44487, 0, 0, 0),
    @(44488, 0, 0, 0),
    @(44489, 0, 0, 0),
    @(44490, 0, 0, 0),
    @(44491, 0, 0, 0),
    @(44492, 0, 0, 0),
    @(44493, 0, 0, 0),
    @(44494, 0, 0, 0),
    @(44495, 0, 0, 0),
    @(44496, 0, 0, 0),
    @(44497, 0, 0, 0),
    @(44498, 1, 1, 15.22997258604935),
    @(44499, 0, 1, 15.22997258604935),
    @(44500, 0, 1, 15.22997258604935),
    @(44501, 0, 1, 15.22997258604935),
    @(44502, 1, 2, 30.45994517209869),
    @(44503, 1, 3, 45.68991775814803),
    @(44504, 2, 5, 76.14986293024673),
    @(44505, 5, 9, 137.0697532744441),
    @(44506, 0, 9, 137.0697532744441),
    @(44507, 11, 20, 304.5994517209869),
    @(44508, 18, 38, 578.7389582698751),
    @(44509, 0, 37, 563.5089856838258),
    @(44510, 0, 36, 548.2790130977764),
    @(44511, 15, 49, 746.2686567164179),
    @(44512, 8, 52, 791.9585744745659),
    @(44513, 6, 58, 883.3384099908621),
    @(44514, 0, 47, 715.8087115443192),
    @(44515, 6, 35, 533.0490405117271),
    @(44516, 4, 39, 593.9689308559244),
    @(44517, 4, 43, 654.8888212001219),
    @(44518, 5, 33, 502.5890953396284),
    @(44519, 7, 32, 487.3591227535791),
    @(44520, 3, 29, 441.6692049954311),
    @(44521, 9, 38, 578.7389582698751),
    @(44522, 11, 43, 654.8888212001219),
    @(44523, 0, 39, 593.9689308559244),
    @(44524, 11, 46, 700.5787389582699),
    @(44525, 2, 43, 654.8888212001219),
    @(44526, 1, 37, 563.5089856838258),
    @(44527, 2, 36, 548.2790130977764),
    @(44528, 0, 27, 411.2092598233323),
    @(44529, 4, 20, 304.5994517209869),
    @(44530, 0, 20, 304.5994517209869),
    @(44531, 0, 9, 137.0697532744441),
    @(44532, 0, 7, 106.6098081023454),
    @(44533, 1, 7, 106.6098081023454),
    @(44534, 1, 6, 91.37983551629607),
    @(44535, 0, 6, 91.37983551629607),
    @(44536, 4, 6, 91.37983551629607),
    @(44537, 1, 7, 106.6098081023454),
    @(44538, 1, 8, 121.8397806883948)
)

$startRow = 386
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $values[0]
    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
}

# Match existing formatting: column A uses the bordered/bold/centered date style
$lastOldRow = $startRow - 1
$lastNewRow = $startRow + $newRows.Count - 1
$ws.Range("A" + $lastOldRow).Copy()
$ws.Range("A" + $startRow + ":A" + $lastNewRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false
